$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = -1.423328640925922
$ws.Range("B2").Value = -1.19501468705647

$ws.Range("A3").Value = -0.5042364550831944
$ws.Range("B3").Value = -0.6391314888068841

$ws.Range("A4").Value = -1.008619983838063
$ws.Range("B4").Value = -0.8452853894249533

$ws.Range("A5").Value = -0.6793797624891781
$ws.Range("B5").Value = -0.6212515765769069

$ws.Range("A6").Value = 0.8175056764132963
$ws.Range("B6").Value = 0.7172543454243656

$ws.Range("A7").Value = -0.05953036942819206
$ws.Range("B7").Value = 0.03269886473398449

$ws.Range("A8").Value = 0.8070948150488942
$ws.Range("B8").Value = 0.6651147549566339

$ws.Range("A9").Value = 0.3089235724971614
$ws.Range("B9").Value = 0.3071880275403656

$ws.Range("A10").Value = -0.1776918520982292
$ws.Range("B10").Value = -0.04165652221259124

$ws.Range("A11").Value = -0.3443318841345794
$ws.Range("B11").Value = -0.2153747805415667

$ws.Range("A12").Value = -0.5091826595782183
$ws.Range("B12").Value = -0.2758957494173607

$ws.Range("A13").Value = 0.426041208712112
$ws.Range("B13").Value = 0.3525548721630504
